# Generate Report for Archive
#
# 1) The "Status" column value changes from "Ready for handoff" to
#    "In Translation" on every sheet/cell that currently shows it.
# 2) The "Status" column (zh-cn / de-de sheets, column C) and the
#    corresponding locale columns on the Overview sheet (columns E/F)
#    get narrower - the stored column width shrinks from
#    17.2159881591797 to 13.4101845877511 "units".
#
# NOTE: Excel's COM ColumnWidth property is expressed in "characters"
# and gets snapped onto the workbook's pixel grid when it is written
# back to the OOXML <col width="..."> attribute (width_xml =
# (Round(chars * MDW) + 5) / MDW, with MDW = 6 for this workbook's
# default font). 13.4101845877511 is not itself representable on that
# grid (it was produced by the original, non-Excel tool that first
# wrote this workbook), so we ask for the "characters" value that is
# mathematically equivalent (target - 5/MDW) and let the engine snap
# it to the nearest grid point for us - that lands on the closest
# achievable width.

$wb = $excel.ActiveWorkbook

$newStatus = "In Translation"
$targetWidthUnits = 13.4101845877511
$mdw = 6
$newColumnWidthChars = $targetWidthUnits - (5 / $mdw)

# --- Overview sheet: locale columns (zh-cn / de-de) show the status too ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidthChars
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidthChars

# --- zh-cn sheet: Status column is column C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidthChars

# --- de-de sheet: Status column is column C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidthChars
